$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(4,19,5,1),
    @(6,8,5,12),
    @(7,5,5,15),
    @(4,14,3,6),
    @(3,12,6,8),
    @(5,2,4,18),
    @(5,6,4,14),
    @(3,6,4,14),
    @(5,5,4,15),
    @(2,7,4,13),
    @(7,7,9,13),
    @(5,4,1,16),
    @(5,6,6,14),
    @(6,5,7,15),
    @(5,14,4,6),
    @(4,8,2,12),
    @(3,13,4,7),
    @(6,14,7,6),
    @(4,12,5,8),
    @(5,5,4,15),
    @(4,7,5,13),
    @(6,6,7,14),
    @(5,7,3,13),
    @(3,4,5,16),
    @(6,6,5,14),
    @(4,7,5,13)
)

$startRow = 1987
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Range("A" + $row).Value = $vals[0]
    $ws.Range("B" + $row).Value = $vals[1]
    $ws.Range("C" + $row).Value = $vals[2]
    $ws.Range("D" + $row).Value = $vals[3]
}

$ws.Range("A2013").Select()
